$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of sorting data for bat 2289 across 3 new dates (3 TTs each, 4 rows per date)
$dates = @(43266, 43270, 43271)
$comments = @("one isoalted cluster", "some MU", "many spikes but hard to isolate")

$dateStyleSrc = $ws.Cells.Item(140, 2)

$row = 141
for ($d = 0; $d -lt $dates.Length; $d++) {
    $dateVal = $dates[$d]
    for ($tt = 1; $tt -le 4; $tt++) {
        $ws.Cells.Item($row, 1).Value = 2289
        $c = $ws.Cells.Item($row, 2)
        $dateStyleSrc.Copy($c)
        $c.Value = $dateVal
        $ws.Cells.Item($row, 3).Value = $tt
        $ws.Cells.Item($row, 4).Value = 6
        $ws.Cells.Item($row, 5).Value = 0
        $ws.Cells.Item($row, 6).Value = 0.8
        $ws.Cells.Item($row, 7).Value = 24
        $ws.Cells.Item($row, 8).Value = 5
        $ws.Cells.Item($row, 9).Value = 4
        $ws.Cells.Item($row, 10).Value = "highpass"
        if ($tt -eq 4) {
            $ws.Cells.Item($row, 11).Value = $comments[$d]
        } else {
            $ws.Cells.Item($row, 11).Value = "no cells"
        }
        $row = $row + 1
    }
}

$win = $excel.ActiveWindow
$win.ScrollRow = 123
$ws.Range("F160").Select()
